# Insert a new data row into the "Sandia" sheet at row 104.
# This shifts the existing rows 104:188 down to 105:189 and keeps their
# data intact. The newly opened row 104 is then filled with a new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 104 (existing row 104 and everything below
# moves down to row 105, etc.)
$ws.Rows("104:104").Insert()

# Populate the new row 104 with the new record's data.
$ws.Range("A104").Value = 11
$ws.Range("B104").Value = "Vega Monumental Concepción"
$ws.Range("C104").Value = "Bíobío"
$ws.Range("D104").Value = 45001
$ws.Range("E104").Value = 8
$ws.Range("F104").Value = 100112028
$ws.Range("G104").Value = "Sandia"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 900
$ws.Range("K104").Value = 1800
$ws.Range("L104").Value = 2000
$ws.Range("M104").Value = 1911
$ws.Range("N104").Value = "$/unidad"
$ws.Range("O104").Value = "Región Metropolitana"
$ws.Range("P104").Value = 1911
$ws.Range("Q104").Value = 1
$ws.Range("R104").Value = "Hortaliza"
